$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = 0.05      # Total P&L %
$summary.Range("B6").Value = 91        # Total Trades
$summary.Range("B9").Value = 49.45     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D6").Value = 30         # Trades
$status.Range("G6").Value = 50         # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet - close trade #91 (row 92) and append the new trade (row 121)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G92").Value = 0.98
$allTrades.Range("H92").Value = "CLOSED"
$allTrades.Range("K92").Value = 99.41
$allTrades.Range("L92").Value = "early_exit"
$allTrades.Range("M92").Value = 0.12

$allTrades.Range("A121").Value = 120
$allTrades.Range("B121").NumberFormat = "@"
$allTrades.Range("B121").Value = "2026-02-18"
$allTrades.Range("B121").ClearFormats()
$allTrades.Range("C121").Value = "00:23:33"
$allTrades.Range("D121").Value = "MarketMaking"
$allTrades.Range("E121").Value = "DOWN"
$allTrades.Range("F121").Value = 0.98
$allTrades.Range("H121").Value = "OPEN"
$allTrades.Range("I121").Value = 0
$allTrades.Range("J121").Value = 0
$allTrades.Range("K121").Value = 99.410254715139
$allTrades.Range("M121").Value = 0
$allTrades.Range("N121").Value = 0
$allTrades.Range("O121").Value = 0
$allTrades.Range("P121").Value = 0.6
$allTrades.Range("Q121").Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet - close trade #91 (row 31) and append the new trade (row 41)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("G31").Value = 0.98
$mm.Range("H31").Value = "CLOSED"
$mm.Range("K31").Value = 99.41
$mm.Range("P31").Value = "early_exit"
$mm.Range("Q31").Value = 0.12

$mm.Range("A41").Value = 120
$mm.Range("B41").NumberFormat = "@"
$mm.Range("B41").Value = "2026-02-18"
$mm.Range("B41").ClearFormats()
$mm.Range("C41").Value = "00:23:33"
$mm.Range("D41").Value = "MarketMaking"
$mm.Range("E41").Value = "DOWN"
$mm.Range("F41").Value = 0.98
$mm.Range("H41").Value = "OPEN"
$mm.Range("I41").Value = 0
$mm.Range("J41").Value = 0
$mm.Range("K41").Value = 99.410254715139
$mm.Range("L41").Value = 0
$mm.Range("M41").Value = 0
$mm.Range("N41").Value = 0.6
$mm.Range("O41").Value = "Normal spread capture: 198 bps"
$mm.Range("Q41").Value = 0
